$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Nora Gavrea): PartIII End date bumped forward ---
$ws.Range("E2").Value = 45444

# --- Row 6 (Charlotte Priestley): PartIII End formula capped at a fixed date ---
$ws.Range("E6").Formula = "=MIN(NOW(), ""01/06/2025  00:00:00"")"

# --- Row 10 (Patrick Whitman): PartIII End formula capped at a fixed date ---
$ws.Range("E10").Formula = "=MIN(NOW(), ""01/06/2025  00:00:00"")"

# --- Row 11 (Lehan Li): PartIII End formula capped at a fixed date ---
$ws.Range("E11").Formula = "=MIN(NOW(), ""01/06/2025  00:00:00"")"

# --- New row 12: Adele Chu, Mphil Start/End, Cosupervision style ---
$ws.Range("A12").Value = "Adele Chu"

$ws.Range("F12").NumberFormat = "mmm-yy"
$ws.Range("F12").Value = 45658

$ws.Range("G12").NumberFormat = "m/d/yy h:mm"
$ws.Range("G12").Formula = "=NOW()"

$ws.Range("J12").Value = "Cosupervision"

# --- Update the active selection to match the new state ---
$ws.Range("E3").Select()
